$d = $word.ActiveDocument

# New, consolidated Catalan text for the "Hercules" campaign dates.
$rsquo = [char]0x2019
$newText = "Dates de la campanya Hercules: 13-22 de juny, 12-21 de juliol, del 10 al 19 d" + $rsquo + "agost"

# Collect the paragraphs whose text starts with "Dates de la campanya"
# (ignoring an optional leading space run present in the very first
# occurrence). This avoids the mid-sentence lowercase occurrence of the
# phrase found elsewhere in the document.
$targets = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Dates de la campanya") -or $t.StartsWith(" Dates de la campanya")) {
        $targets += $p
    }
}

foreach ($p in $targets) {
    $rng = $p.Range
    # Exclude the trailing paragraph-mark character so only the runs are
    # replaced, keeping the paragraph's own <w:pPr> untouched.
    $inner = $d.Range($rng.Start, $rng.End - 1)

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $inner.InsertXML($xmlFrag)
}
